$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row 2 and row 3 for the affected columns
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)

$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
